$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "as at" date in the intro paragraph (A2)
$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 17 April 2025"

# 2. Insert a new row at row 37 (pushes old rows 37..64 down to 38..65),
#    matching the style/format of the surrounding rows.
$ws.Range("A37:F37").EntireRow.Insert()

# 3. Populate the newly inserted row 37 with the new publication entry,
#    matching its sibling row (row 36) for the week-commencing / publication
#    date / status / type columns.
$ws.Range("A37").Value = "28 Jul 2025"
$ws.Range("B37").Value = "HMPPS Annual Digest, April 2024 to March 2025"
$ws.Range("C37").Value = "31 July 2025"
$ws.Range("D37").Value = "provisional"
$ws.Range("E37").Value = 31
$ws.Range("F37").Value = "standard"

# 4. Extend the conditional formatting ranges to cover the new last row (65)
#    instead of the old last row (64).
$fc = $ws.Range("A5:F64").FormatConditions
$fc.Item(1).ModifyAppliesToRange($ws.Range("A5:F65"))
$fc.Item(2).ModifyAppliesToRange($ws.Range("A5:F65"))
$fc.Item(3).ModifyAppliesToRange($ws.Range("A5:F65"))
$fc.Item(4).ModifyAppliesToRange($ws.Range("A5:A65"))
$fc.Item(5).ModifyAppliesToRange($ws.Range("A5:A65"))
